$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1939252336448598
$ws.Range("C2").Value = 0.5911214953271028
$ws.Range("J2").Value = 0.01635514018691589
$ws.Range("P2").Value = 0.1331775700934579
$ws.Range("S2").Value = 0.06542056074766354
$ws.Range("C3").Value = 0.01886792452830189
$ws.Range("J3").Value = 0.04150943396226415
$ws.Range("P3").Value = 0.7320754716981132
$ws.Range("S3").Value = 0.2075471698113208
$ws.Range("J4").Value = 0.02127659574468085
$ws.Range("P4").Value = 0.6382978723404256
$ws.Range("S4").Value = 0.3404255319148936
$ws.Range("B6").Value = 0.07224334600760456
$ws.Range("D6").Value = 0.01140684410646388
$ws.Range("F6").Value = 0.06463878326996197
$ws.Range("J6").Value = 0.3155893536121673
$ws.Range("O6").Value = 0.01901140684410646
$ws.Range("Q6").Value = 0.1178707224334601
$ws.Range("R6").Value = 0.05703422053231939
$ws.Range("S6").Value = 0.3422053231939163
$ws.Range("B7").Value = 0.1487603305785124
$ws.Range("D7").Value = 0.008264462809917356
$ws.Range("F7").Value = 0.008264462809917356
$ws.Range("J7").Value = 0.1446280991735537
$ws.Range("O7").Value = 0.02066115702479339
$ws.Range("Q7").Value = 0.115702479338843
$ws.Range("R7").Value = 0.08264462809917356
$ws.Range("S7").Value = 0.4710743801652892
$ws.Range("B8").Value = 0.1276102088167053
$ws.Range("D8").Value = 0.02784222737819025
$ws.Range("E8").Value = 0.006960556844547564
$ws.Range("F8").Value = 0.04176334106728538
$ws.Range("J8").Value = 0.1276102088167053
$ws.Range("O8").Value = 0.01624129930394431
$ws.Range("Q8").Value = 0.1902552204176334
$ws.Range("R8").Value = 0.05568445475638051
$ws.Range("S8").Value = 0.4060324825986079
$ws.Range("B9").Value = 0.1510416666666667
$ws.Range("D9").Value = 0.01041666666666667
$ws.Range("F9").Value = 0.05208333333333334
$ws.Range("J9").Value = 0.109375
$ws.Range("O9").Value = 0.02083333333333333
$ws.Range("Q9").Value = 0.2083333333333333
$ws.Range("R9").Value = 0.109375
$ws.Range("S9").Value = 0.3385416666666667
$ws.Range("B10").Value = 0.1439446366782007
$ws.Range("D10").Value = 0.02076124567474048
$ws.Range("F10").Value = 0.07750865051903114
$ws.Range("J10").Value = 0.1038062283737024
$ws.Range("O10").Value = 0.01453287197231834
$ws.Range("Q10").Value = 0.1889273356401384
$ws.Range("R10").Value = 0.08027681660899653
$ws.Range("S10").Value = 0.370242214532872
$ws.Range("G11").Value = 0.1436619718309859
$ws.Range("J11").Value = 0.07042253521126761
$ws.Range("K11").Value = 0.1887323943661972
$ws.Range("L11").Value = 0.5859154929577465
$ws.Range("S11").Value = 0.01126760563380282
$ws.Range("G12").Value = 0.7674418604651163
$ws.Range("J12").Value = 0.1813953488372093
$ws.Range("K12").Value = 0.004651162790697674
$ws.Range("L12").Value = 0.02790697674418605
$ws.Range("S12").Value = 0.0186046511627907
$ws.Range("G13").Value = 0.7021276595744681
$ws.Range("J13").Value = 0.2765957446808511
$ws.Range("S13").Value = 0.02127659574468085
$ws.Range("F15").Value = 0.01639344262295082
$ws.Range("H15").Value = 0.09836065573770492
$ws.Range("I15").Value = 0.06147540983606557
$ws.Range("J15").Value = 0.4098360655737705
$ws.Range("K15").Value = 0.08196721311475409
$ws.Range("M15").Value = 0.00819672131147541
$ws.Range("N15").Value = 0.004098360655737705
$ws.Range("O15").Value = 0.09836065573770492
$ws.Range("S15").Value = 0.2213114754098361
$ws.Range("F16").Value = 0.02205882352941177
$ws.Range("H16").Value = 0.1654411764705882
$ws.Range("I16").Value = 0.05882352941176471
$ws.Range("J16").Value = 0.4264705882352941
$ws.Range("K16").Value = 0.1102941176470588
$ws.Range("M16").Value = 0.01470588235294118
$ws.Range("O16").Value = 0.05882352941176471
$ws.Range("S16").Value = 0.1433823529411765
$ws.Range("F17").Value = 0.03070175438596491
$ws.Range("H17").Value = 0.162280701754386
$ws.Range("I17").Value = 0.09210526315789473
$ws.Range("J17").Value = 0.4254385964912281
$ws.Range("K17").Value = 0.09429824561403509
$ws.Range("M17").Value = 0.02412280701754386
$ws.Range("N17").Value = 0.002192982456140351
$ws.Range("O17").Value = 0.04824561403508772
$ws.Range("S17").Value = 0.1206140350877193
$ws.Range("F18").Value = 0.02564102564102564
$ws.Range("H18").Value = 0.158974358974359
$ws.Range("I18").Value = 0.07179487179487179
$ws.Range("J18").Value = 0.4717948717948718
$ws.Range("K18").Value = 0.1025641025641026
$ws.Range("M18").Value = 0.01538461538461539
$ws.Range("O18").Value = 0.06666666666666667
$ws.Range("S18").Value = 0.08717948717948718
$ws.Range("F19").Value = 0.03012912482065997
$ws.Range("H19").Value = 0.1886657101865136
$ws.Range("I19").Value = 0.07747489239598278
$ws.Range("J19").Value = 0.3694404591104735
$ws.Range("K19").Value = 0.1226685796269727
$ws.Range("M19").Value = 0.02008608321377331
$ws.Range("N19").Value = 0.001434720229555237
$ws.Range("O19").Value = 0.07532281205164992
$ws.Range("S19").Value = 0.1147776183644189
